$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 8 (shifting old rows 8-10 down to 10-12)
$ws.Rows("8:9").Insert()

# Fill new rows column by column (A8, A9, B8, B9, C8, C9, D8, D9)
$ws.Range("A8").Value = "102_AutomobileInsurance_004_ProductData_001_MandatoryFields"
$ws.Range("A9").Value = "102_AutomobileInsurance_004_ProductData_002_FieldHintsAndErrors"

$ws.Range("B8").Value = "var102_AutomobileInsurance_004_ProductData_001_MandatoryFields"
$ws.Range("B9").Value = "var102_AutomobileInsurance_004_ProductData_002_FieldHintsAndErrors"

$ws.Range("C8").Value = "Open Automobile Insurance"
$ws.Range("C9").Value = "Open Automobile Insurance"

$ws.Range("D8").Value = "102_AutomobileInsurance_004_ProductData_001_MandatoryFields"
$ws.Range("D9").Value = "102_AutomobileInsurance_004_ProductData_002_FieldHintsAndErrors"

# Set selection to A9 (matches the saved selection in the target workbook)
$ws.Range("A9").Select()
